# Update the "Overview" worksheet values with the latest financial data
# (adds the monte_carlo period column values / updates the yearly rial database)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Row 11: فروش (Sales)
$ws.Range("D11").Value = 9193650
$ws.Range("E11").Value = 12077819
$ws.Range("F11").Value = 16860698
$ws.Range("G11").Value = 33645124
$ws.Range("H11").Value = 49765416

# Row 12: بهای تمام شده کالای فروش رفته (Cost of goods sold)
$ws.Range("D12").Value = -7455830
$ws.Range("E12").Value = -10240368
$ws.Range("F12").Value = -14082370
$ws.Range("G12").Value = -26395550
$ws.Range("H12").Value = -46833233

# Row 13: سود (زیان) ناخالص (Gross profit)
$ws.Range("D13").Value = 1737820
$ws.Range("E13").Value = 1837451
$ws.Range("F13").Value = 2778328
$ws.Range("G13").Value = 7249574
$ws.Range("H13").Value = 2932183

# Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses)
$ws.Range("D14").Value = -331419
$ws.Range("E14").Value = -417810
$ws.Range("F14").Value = -571114
$ws.Range("G14").Value = -919142
$ws.Range("H14").Value = -1251487

# Row 15: هزینه کاهش ارزش دریافتنی‌‏ها (هزینه استثنایی) - E:H become numeric 0, D stays "-"
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0

# Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی
$ws.Range("D16").Value = -6424
$ws.Range("E16").Value = -761279
$ws.Range("F16").Value = -73980
$ws.Range("G16").Value = 56470
$ws.Range("H16").Value = -55475

# Row 17: سود (زیان) عملیاتی (Operating profit)
$ws.Range("D17").Value = 1399977
$ws.Range("E17").Value = 658362
$ws.Range("F17").Value = 2133234
$ws.Range("G17").Value = 6386902
$ws.Range("H17").Value = 1625221

# Row 18: هزینه های مالی (Financial expenses)
$ws.Range("D18").Value = -148269
$ws.Range("E18").Value = -243317
$ws.Range("F18").Value = -406797
$ws.Range("G18").Value = -696912
$ws.Range("H18").Value = -1610584

# Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی
$ws.Range("D19").Value = 44841
$ws.Range("E19").Value = 386533
$ws.Range("F19").Value = 237290
$ws.Range("G19").Value = 1989668
$ws.Range("H19").Value = 1447920

# Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات
$ws.Range("D20").Value = 1296549
$ws.Range("E20").Value = 801578
$ws.Range("F20").Value = 1963727
$ws.Range("G20").Value = 7679658
$ws.Range("H20").Value = 1462557

# Row 21: مالیات (Tax)
$ws.Range("D21").Value = -313882
$ws.Range("E21").Value = -101152
$ws.Range("F21").Value = -255918
$ws.Range("G21").Value = -1069506
$ws.Range("H21").Value = -39937

# Row 22: سود (زیان) خالص عملیات در حال تداوم
$ws.Range("D22").Value = 982667
$ws.Range("E22").Value = 700426
$ws.Range("F22").Value = 1707809
$ws.Range("G22").Value = 6610152
$ws.Range("H22").Value = 1422620

# Row 23: سود (زیان) عملیات متوقف شده پس از اثر مالیاتی - all columns become numeric 0
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0

# Row 24: سود (زیان) خالص (Net profit)
$ws.Range("D24").Value = 982667
$ws.Range("E24").Value = 700426
$ws.Range("F24").Value = 1707809
$ws.Range("G24").Value = 6610152
$ws.Range("H24").Value = 1422620

# Row 25: سود هر سهم پس از کسر مالیات
$ws.Range("D25").Value = 389
$ws.Range("E25").Value = 277
$ws.Range("F25").Value = 676
$ws.Range("G25").Value = 2616
$ws.Range("H25").Value = 563

# Row 26: سرمایه (Capital)
$ws.Range("D26").Value = 2526500
$ws.Range("E26").Value = 2526500
$ws.Range("F26").Value = 2526500
$ws.Range("G26").Value = 2526500
$ws.Range("H26").Value = 2526500

# Row 27: سود هر سهم بر اساس آخرین سرمایه
$ws.Range("D27").Value = 389
$ws.Range("E27").Value = 277
$ws.Range("F27").Value = 676
$ws.Range("G27").Value = 2616
$ws.Range("H27").Value = 563
